$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.354.25'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.940.51'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.08%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7168'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -7.34%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3342'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.67'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07358'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8157'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08152'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.940.30'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.493'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.31'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.89'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.379.00'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008376'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -7.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.858'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.196.39'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.977'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.852'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.54'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.410'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.40'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1315'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -10.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.576'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.345'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.490'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.268'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05272'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.278'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7624'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.758'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01996'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.848'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '81.35'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.563'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4568'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.034'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8471'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.905'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.468'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.16'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4190'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.76%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06053'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.35%  '
